$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Five row-pairs had their match data (columns F..V) swapped between the
#    two rows (the "Indice"/row-position in column A, and B..E, stay put).
# ---------------------------------------------------------------------------
$swapPairs = @(
    @(3, 4),
    @(15, 16),
    @(28, 30),
    @(42, 44),
    @(46, 48)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 6; $col -le 22; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# ---------------------------------------------------------------------------
# 2) Seven brand-new match rows were appended at the bottom of the sheet
#    (rows 55..61), extending the used range from A1:V54 to A1:V61.
# ---------------------------------------------------------------------------
$newRows = @(
    @(54, "rwanda", "premier-league", "2023-2024", 45226.625, "Gorilla", 2, "Marines", 3, 2.33, "26/10/2023 03:12", 2.49, "27/10/2023 14:39", 2.81, "26/10/2023 03:12", 2.7, "27/10/2023 14:46", 2.77, "26/10/2023 03:12", 2.95, "27/10/2023 14:45", "https://www.betexplorer.com/football/rwanda/premier-league/gorilla-marines/rFldkuqE/"),
    @(55, "rwanda", "premier-league", "2023-2024", 45226.75, "Kiyovu", 6, "Etoile de L'Est", 1, 1.49, "26/10/2023 06:12", 1.66, "27/10/2023 17:32", 3.47, "26/10/2023 06:12", 3.25, "27/10/2023 17:32", 5.25, "26/10/2023 06:12", 5.01, "27/10/2023 17:32", "https://www.betexplorer.com/football/rwanda/premier-league/kiyovu-etoile-de-l-est/dIh0labK/"),
    @(56, "rwanda", "premier-league", "2023-2024", 45227.625, "Bugesera", 2, "Amagaju", 2, 1.88, "27/10/2023 03:12", 1.85, "28/10/2023 14:10", 2.89, "27/10/2023 03:12", 3.06, "28/10/2023 14:10", 3.69, "27/10/2023 03:12", 4.16, "28/10/2023 14:10", "https://www.betexplorer.com/football/rwanda/premier-league/bugesera-amagaju/juh4mJDQ/"),
    @(57, "rwanda", "premier-league", "2023-2024", 45227.625, "Mukura Victory Sports", 1, "Etincelles", 0, 1.7, "27/10/2023 03:12", 1.65, "28/10/2023 12:08", 3.15, "27/10/2023 03:12", 3.27, "28/10/2023 13:02", 4.1, "27/10/2023 03:12", 5.02, "28/10/2023 12:08", "https://www.betexplorer.com/football/rwanda/premier-league/mukura-victory-sports-etincelles/hlvWrwMs/"),
    @(58, "rwanda", "premier-league", "2023-2024", 45228.58333333334, "APR", 0, "Rayon Sport", 0, 1.83, "28/10/2023 03:12", 2.06, "29/10/2023 13:56", 2.97, "28/10/2023 03:12", 2.92, "29/10/2023 13:56", 3.77, "28/10/2023 03:12", 3.62, "29/10/2023 13:56", "https://www.betexplorer.com/football/rwanda/premier-league/apr-rayon-sport/b19WtHjf/"),
    @(59, "rwanda", "premier-league", "2023-2024", 45228.58333333334, "Musanze", 1, "AS Kigali", 0, 2.78, "28/10/2023 03:12", 2.97, "29/10/2023 13:46", 2.6, "28/10/2023 03:12", 2.53, "29/10/2023 13:46", 2.54, "28/10/2023 03:12", 2.69, "29/10/2023 13:46", "https://www.betexplorer.com/football/rwanda/premier-league/musanze-as-kigali/6RCvueL6/"),
    @(60, "rwanda", "premier-league", "2023-2024", 45228.58333333334, "Sunrise", 1, "Muhazi United", 2, 2.01, "28/10/2023 03:12", 2.05, "29/10/2023 04:30", 2.82, "28/10/2023 03:12", 2.92, "29/10/2023 12:02", 3.39, "28/10/2023 03:12", 3.59, "29/10/2023 04:30", "https://www.betexplorer.com/football/rwanda/premier-league/sunrise-muhazi-united/hO8zty60/")
)

$destRow = 55
foreach ($rowData in $newRows) {
    # Clone formatting (styles, number formats, borders, ...) from an
    # existing data row so the new row matches the sheet's look & feel.
    $ws.Range("A2:V2").Copy($ws.Range("A" + $destRow + ":V" + $destRow))

    for ($i = 0; $i -lt $rowData.Length; $i++) {
        $ws.Cells.Item($destRow, $i + 1).Value = $rowData[$i]
    }

    $destRow++
}
